$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: financial period headers (shift left one period, add newest period) ---
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (shift left one period, add newest date) ---
$ws.Range("D9").Value = "1400-10-29 (3)"
$ws.Range("E9").Value = "1401-03-11 (8)"
$ws.Range("F9").Value = "1401-04-29 (2)"
$ws.Range("G9").Value = "1401-08-29 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-30 (7)"
$ws.Range("J9").Value = "1401-04-29"
$ws.Range("K9").Value = "1401-08-29 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-30"

# --- Data rows 11-27: shift all quarterly figures left by one column, append newest-quarter figures ---
$ws.Range("D11").Value = 5089178 ; $ws.Range("E11").Value = 8468536 ; $ws.Range("F11").Value = 2522162 ; $ws.Range("G11").Value = 5873558 ; $ws.Range("H11").Value = 9590983 ; $ws.Range("I11").Value = 13691281 ; $ws.Range("J11").Value = 4321632 ; $ws.Range("K11").Value = 10529704 ; $ws.Range("L11").Value = 17172835 ; $ws.Range("M11").Value = 24645941
$ws.Range("D12").Value = -4756469 ; $ws.Range("E12").Value = -7313605 ; $ws.Range("F12").Value = -2139565 ; $ws.Range("G12").Value = -4901473 ; $ws.Range("H12").Value = -8127904 ; $ws.Range("I12").Value = -11319270 ; $ws.Range("J12").Value = -3844001 ; $ws.Range("K12").Value = -8108088 ; $ws.Range("L12").Value = -13797526 ; $ws.Range("M12").Value = -19785175
$ws.Range("D13").Value = 332709 ; $ws.Range("E13").Value = 1154931 ; $ws.Range("F13").Value = 382597 ; $ws.Range("G13").Value = 972085 ; $ws.Range("H13").Value = 1463079 ; $ws.Range("I13").Value = 2372011 ; $ws.Range("J13").Value = 477631 ; $ws.Range("K13").Value = 2421616 ; $ws.Range("L13").Value = 3375309 ; $ws.Range("M13").Value = 4860766
$ws.Range("D14").Value = -419558 ; $ws.Range("E14").Value = -791950 ; $ws.Range("F14").Value = -196255 ; $ws.Range("G14").Value = -471374 ; $ws.Range("H14").Value = -732002 ; $ws.Range("I14").Value = -1093286 ; $ws.Range("J14").Value = -257142 ; $ws.Range("K14").Value = -675717 ; $ws.Range("L14").Value = -982990 ; $ws.Range("M14").Value = -1359815
$ws.Range("D15").Value = 0 ; $ws.Range("E15").Value = 0 ; $ws.Range("F15").Value = 0 ; $ws.Range("G15").Value = 0 ; $ws.Range("H15").Value = 0 ; $ws.Range("I15").Value = 0 ; $ws.Range("J15").Value = 0 ; $ws.Range("K15").Value = 0 ; $ws.Range("L15").Value = 0 ; $ws.Range("M15").Value = 0
$ws.Range("D16").Value = 458206 ; $ws.Range("E16").Value = 411446 ; $ws.Range("F16").Value = -136295 ; $ws.Range("G16").Value = 2773 ; $ws.Range("H16").Value = 34932 ; $ws.Range("I16").Value = 4537 ; $ws.Range("J16").Value = 22552 ; $ws.Range("K16").Value = 35910 ; $ws.Range("L16").Value = 303979 ; $ws.Range("M16").Value = 883120
$ws.Range("D17").Value = 371357 ; $ws.Range("E17").Value = 774427 ; $ws.Range("F17").Value = 50047 ; $ws.Range("G17").Value = 503484 ; $ws.Range("H17").Value = 766009 ; $ws.Range("I17").Value = 1283262 ; $ws.Range("J17").Value = 243041 ; $ws.Range("K17").Value = 1781809 ; $ws.Range("L17").Value = 2696298 ; $ws.Range("M17").Value = 4384071
$ws.Range("D18").Value = -16659 ; $ws.Range("E18").Value = -17293 ; $ws.Range("F18").Value = -10530 ; $ws.Range("G18").Value = -22944 ; $ws.Range("H18").Value = -32702 ; $ws.Range("I18").Value = -48202 ; $ws.Range("J18").Value = -8688 ; $ws.Range("K18").Value = -13629 ; $ws.Range("L18").Value = -15642 ; $ws.Range("M18").Value = -15641
$ws.Range("D19").Value = 109038 ; $ws.Range("E19").Value = 148974 ; $ws.Range("F19").Value = 72750 ; $ws.Range("G19").Value = 81857 ; $ws.Range("H19").Value = 132974 ; $ws.Range("I19").Value = 215534 ; $ws.Range("J19").Value = 27175 ; $ws.Range("K19").Value = 53899 ; $ws.Range("L19").Value = 127306 ; $ws.Range("M19").Value = 141080
$ws.Range("D20").Value = 463736 ; $ws.Range("E20").Value = 906108 ; $ws.Range("F20").Value = 112267 ; $ws.Range("G20").Value = 562397 ; $ws.Range("H20").Value = 866281 ; $ws.Range("I20").Value = 1450594 ; $ws.Range("J20").Value = 261528 ; $ws.Range("K20").Value = 1822079 ; $ws.Range("L20").Value = 2807962 ; $ws.Range("M20").Value = 4509510
$ws.Range("D21").Value = 0 ; $ws.Range("E21").Value = -5083 ; $ws.Range("F21").Value = 0 ; $ws.Range("G21").Value = -32534 ; $ws.Range("H21").Value = -66472 ; $ws.Range("I21").Value = -92886 ; $ws.Range("J21").Value = -33294 ; $ws.Range("K21").Value = -262879 ; $ws.Range("L21").Value = -372643 ; $ws.Range("M21").Value = -510870
$ws.Range("D22").Value = 463736 ; $ws.Range("E22").Value = 901025 ; $ws.Range("F22").Value = 112267 ; $ws.Range("G22").Value = 529863 ; $ws.Range("H22").Value = 799809 ; $ws.Range("I22").Value = 1357708 ; $ws.Range("J22").Value = 228234 ; $ws.Range("K22").Value = 1559200 ; $ws.Range("L22").Value = 2435319 ; $ws.Range("M22").Value = 3998640
$ws.Range("D23").Value = 0 ; $ws.Range("E23").Value = 0 ; $ws.Range("F23").Value = 0 ; $ws.Range("G23").Value = 0 ; $ws.Range("H23").Value = 0 ; $ws.Range("I23").Value = 0 ; $ws.Range("J23").Value = 0 ; $ws.Range("K23").Value = 0 ; $ws.Range("L23").Value = 0 ; $ws.Range("M23").Value = 0
$ws.Range("D24").Value = 463736 ; $ws.Range("E24").Value = 901025 ; $ws.Range("F24").Value = 112267 ; $ws.Range("G24").Value = 529863 ; $ws.Range("H24").Value = 799809 ; $ws.Range("I24").Value = 1357708 ; $ws.Range("J24").Value = 228234 ; $ws.Range("K24").Value = 1559200 ; $ws.Range("L24").Value = 2435319 ; $ws.Range("M24").Value = 3998640
$ws.Range("D25").Value = 824 ; $ws.Range("E25").Value = 493 ; $ws.Range("F25").Value = 61 ; $ws.Range("G25").Value = 290 ; $ws.Range("H25").Value = 438 ; $ws.Range("I25").Value = 744 ; $ws.Range("J25").Value = 125 ; $ws.Range("K25").Value = 854 ; $ws.Range("L25").Value = 1334 ; $ws.Range("M25").Value = 2190
$ws.Range("D26").Value = 563000 ; $ws.Range("E26").Value = 1826000 ; $ws.Range("F26").Value = 1826000 ; $ws.Range("G26").Value = 1826000 ; $ws.Range("H26").Value = 1826000 ; $ws.Range("I26").Value = 1826000 ; $ws.Range("J26").Value = 1826000 ; $ws.Range("K26").Value = 1826000 ; $ws.Range("L26").Value = 1826000 ; $ws.Range("M26").Value = 1826000
$ws.Range("D27").Value = 254 ; $ws.Range("E27").Value = 493 ; $ws.Range("F27").Value = 61 ; $ws.Range("G27").Value = 290 ; $ws.Range("H27").Value = 438 ; $ws.Range("I27").Value = 744 ; $ws.Range("J27").Value = 125 ; $ws.Range("K27").Value = 854 ; $ws.Range("L27").Value = 1334 ; $ws.Range("M27").Value = 2190
